$d = $word.ActiveDocument

# The "download a folder from the server" instructions originally read
# (note: both dashes below are EN DASHES, U+2013):
#
#     pscp –scp -r ml0901@cartesius.surfsara.nl:/home/ml0901/mapnaam ...
#
# "–scp" is a leftover mistake (a stray duplicate of the "pscp" fix), and
# the line should simply read:
#
#     pscp -r ml0901@cartesius.surfsara.nl:/home/ml0901/mapnaam ...
#
# A similar-looking "upload" command earlier in the document ("Pscp –scp
# –r mapnaam ...") must be left alone. There, the dash right before "r"
# is also an en dash, whereas the text we need to fix uses a plain
# hyphen before the "r" ("-r"), so searching for the literal
# en dash + "scp -r" sequence uniquely targets only the spot to fix.

$enDash = [char]0x2013
$findText = $enDash + "scp -r"

$d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "-r", 2)
